# Update the "取得日時" (acquired-at) timestamp column on the first sheet
# (ランサーズ) from the previous scrape run to the new one.
#
# Commit: Append: 2025-10-12 01:46 JST
#
# All rows 2-14 in column A previously held "2025-10-12 01:18:19" and are
# being bumped to the new run's timestamp "2025-10-12 01:46:49".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-12 01:18:19"
$newTimestamp = "2025-10-12 01:46:49"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
